# Update automatico via Actualizar 02-06-2021 12-17-13
#
# The "Fecha" (date) column D holds one timestamp per 14-row availability
# check block (rows 2-15, 16-29, 30-43). A fresh check cycle ran: the newest
# block (rows 2-15) gets a brand-new timestamp, and the two older blocks
# shift down to take on the timestamp that used to belong to the block
# above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44233.51190966772
$ws.Range("D16:D29").Value = 44233.49073519676
$ws.Range("D30:D43").Value = 44233.46954513889
